{"js": "// The document contains two references to the pom file naming:\n//   1) \"This is the content of the pom.xml.jim file:\"\n//   2) \"{@include [verbatim] ../../../pom.xml.jam}}\"\n// The edit drops the redundant \".xml\" segment from both file names so that\n// both now read \"pom.jam\" (per the commit message, the jbang/sh template\n// no longer needs the \".xml\" part of the generated file name, since the\n// version is now read from version.jim at the top level).\n\n// --- 1) \"pom.xml.jim\" -> \"pom.jam\" (note: the final three runs read\n//        \"pom.j\" + \"a\" + \"m\", i.e. the visible text becomes \"pom.jam\") ---\nconst jimResults = context.document.body.search(\"pom.xml.jim\", { matchCase: true });\njimResults.load(\"text\");\nawait context.sync();\n\nif (jimResults.items.length > 0) {\n  jimResults.items[0].insertText(\"pom.jam\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) \"pom.xml.jam\" -> \"pom.jam\" ---\nconst jamResults = context.document.body.search(\"pom.xml.jam\", { matchCase: true });\njamResults.load(\"text\");\nawait context.sync();\n\nif (jamResults.items.length > 0) {\n  jamResults.items[0].insertText(\"pom.jam\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document contains two references to the pom file naming:\n#   1) \"This is the content of the pom.xml.jim file:\"\n#   2) \"{@include [verbatim] ../../../pom.xml.jam}}\"\n# The edit drops the redundant \".xml\" segment from both file names so that\n# the visible text in both places reads \"pom.jam\" (per the commit message,\n# the generated file name no longer needs the \".xml\" part).\n\n$d = $word.ActiveDocument\n\n# --- 1) \"pom.xml.jim\" -> \"pom.jam\" ---\n$find1 = $d.Content.Find\n$find1.Text = \"pom.xml.jim\"\n$find1.Replacement.Text = \"pom.jam\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, \"pom.jam\", 2)\n\n# --- 2) \"pom.xml.jam\" -> \"pom.jam\" ---\n$find2 = $d.Content.Find\n$find2.Text = \"pom.xml.jam\"\n$find2.Replacement.Text = \"pom.jam\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, \"pom.jam\", 2)\n"}
